# Add a new "Turkey" worksheet (Zettler test data) after "Spain",
# based on a copy of the "Spain" sheet so formatting/shared styles line up.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Duplicate the Spain sheet and place the copy right after it; Excel
# names the copy "Spain (2)" and makes it the active sheet.
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Turkey-specific content.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3310"

# The shorter strings no longer need the taller wrapped row height that
# Spain's longer text required - restore rows 3-5 to the sheet default.
$turkey.Rows.Item(3).EntireRow.AutoFit()
$turkey.Rows.Item(4).EntireRow.AutoFit()
$turkey.Rows.Item(5).EntireRow.AutoFit()

# Column D narrows slightly to fit the new content (serializes to width=22).
$turkey.Columns.Item(4).ColumnWidth = 21.166666666666668

# Restore Spain's view (no longer the active tab) to a full-range selection,
# and leave Turkey selected at G15 as the active sheet/tab.
$spain.Select()
$spain.Range("A1:D12").Select()

$turkey.Select()
$turkey.Range("G15").Select()
